# Auto-generated edit script applying market-data refresh to Cactuar_Profits sheets
# Values below come from the authoritative diff (old -> new) per cell.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 69995
$ws.Range("I68").Value = 69995
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 69995
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("M68").Value = -69246
$ws.Range("H69").Value = 3997.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3997.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 11992.5
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -13740.5
$ws.Range("H71").Value = 69995
$ws.Range("I71").Value = 69995
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 209985
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("M71").Value = -206241
$ws.Range("H72").Value = 3997.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3997.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 35977.5
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -44713.5
$ws.Range("H112").Value = 3675.6128
$ws.Range("J112").Value = 3723.1667
$ws.Range("L112").Value = 11169.5001
$ws.Range("N112").Value = -13385.5001
$ws.Range("H129").Value = 1712.7778
$ws.Range("J129").Value = 3463.6667
$ws.Range("L129").Value = 10391.0001
$ws.Range("N129").Value = -20391.0001
$ws.Range("H135").Value = 6047.143
$ws.Range("J135").Value = 11219.1
$ws.Range("L135").Value = 100971.9
$ws.Range("N135").Value = -106041.9
$ws.Range("H137").Value = 19612684
$ws.Range("I137").Value = 1833.3334
$ws.Range("K137").Value = 5500.0002
$ws.Range("M137").Value = -2950.0002
$ws.Range("H138").Value = 5685.288
$ws.Range("I138").Value = 1884.9412
$ws.Range("J138").Value = 7223.524
$ws.Range("K138").Value = 5654.8236
$ws.Range("L138").Value = 21670.572
$ws.Range("M138").Value = -514.8235999999997
$ws.Range("N138").Value = -31950.572
$ws.Range("H140").Value = 71468.09
$ws.Range("J140").Value = 68981.11
$ws.Range("L140").Value = 68981.11
$ws.Range("N140").Value = -79341.11
$ws.Range("H141").Value = 6845.4165
$ws.Range("I141").Value = 5905
$ws.Range("K141").Value = 17715
$ws.Range("M141").Value = -12535

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 949252.25
$ws.Range("I2").Value = 1154884.2
$ws.Range("K2").Value = 1154884.2
$ws.Range("M2").Value = -1154771.2
$ws.Range("H32").Value = 3940.209
$ws.Range("I32").Value = 2199.7778
$ws.Range("K32").Value = 2199.7778
$ws.Range("M32").Value = -1912.7778
$ws.Range("H45").Value = 1657.3636
$ws.Range("I45").Value = 1440.875
$ws.Range("K45").Value = 1440.875
$ws.Range("M45").Value = -1063.875
$ws.Range("H74").Value = 28847322
$ws.Range("I74").Value = 37501064
$ws.Range("J74").Value = 1515.6666
$ws.Range("K74").Value = 37501064
$ws.Range("L74").Value = 1515.6666
$ws.Range("M74").Value = -37500190
$ws.Range("N74").Value = -3263.6666
$ws.Range("H77").Value = 28847322
$ws.Range("I77").Value = 37501064
$ws.Range("J77").Value = 1515.6666
$ws.Range("K77").Value = 187505320
$ws.Range("L77").Value = 7578.333000000001
$ws.Range("M77").Value = -187500952
$ws.Range("N77").Value = -16314.333
$ws.Range("H102").Value = 1371664.9
$ws.Range("I102").Value = 1958232.6
$ws.Range("J102").Value = 3007
$ws.Range("K102").Value = 1958232.6
$ws.Range("L102").Value = 3007
$ws.Range("M102").Value = -1956610.6
$ws.Range("N102").Value = -6251
$ws.Range("H116").Value = 949252.25
$ws.Range("I116").Value = 1154884.2
$ws.Range("K116").Value = 1154884.2
$ws.Range("M116").Value = -1152590.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 949252.25
$ws.Range("I3").Value = 1154884.2
$ws.Range("K3").Value = 1154884.2
$ws.Range("M3").Value = -1154770.2
$ws.Range("H105").Value = 1685.8055
$ws.Range("I105").Value = 1208.9333
$ws.Range("K105").Value = 1208.9333
$ws.Range("M105").Value = 538.0667000000001
$ws.Range("H107").Value = 26317448
$ws.Range("I107").Value = 55557856
$ws.Range("J107").Value = 1082.4
$ws.Range("K107").Value = 55557856
$ws.Range("L107").Value = 1082.4
$ws.Range("M107").Value = -55555936
$ws.Range("N107").Value = -4922.4
$ws.Range("H134").Value = 1408.5834
$ws.Range("I134").Value = 1210.5555
$ws.Range("J134").Value = 2002.6666
$ws.Range("K134").Value = 3631.6665
$ws.Range("L134").Value = 6007.9998
$ws.Range("M134").Value = -1096.6665
$ws.Range("N134").Value = -11077.9998
$ws.Range("H140").Value = 168166.75
$ws.Range("J140").Value = 168166.75
$ws.Range("L140").Value = 168166.75
$ws.Range("N140").Value = -178526.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3323.361
$ws.Range("I31").Value = 3191.3872
$ws.Range("J31").Value = 4141.6
$ws.Range("K31").Value = 3191.3872
$ws.Range("L31").Value = 4141.6
$ws.Range("M31").Value = -2896.3872
$ws.Range("N31").Value = -4731.6
$ws.Range("H34").Value = 3323.361
$ws.Range("I34").Value = 3191.3872
$ws.Range("J34").Value = 4141.6
$ws.Range("K34").Value = 3191.3872
$ws.Range("L34").Value = 4141.6
$ws.Range("M34").Value = -2989.3872
$ws.Range("N34").Value = -4545.6
$ws.Range("H58").Value = 2502702.5
$ws.Range("I58").Value = 2502702.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2502702.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2502499.5
$ws.Range("N58").Value = ""
$ws.Range("H105").Value = 4133729.8
$ws.Range("I105").Value = 7576588
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 7576588
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = -7574841
$ws.Range("N105").Value = -5794
$ws.Range("H132").Value = 2024.3636
$ws.Range("I132").Value = 1211.3334
$ws.Range("K132").Value = 3634.0002
$ws.Range("M132").Value = -1104.0002
$ws.Range("H134").Value = 3660.923
$ws.Range("J134").Value = 4374.25
$ws.Range("L134").Value = 13122.75
$ws.Range("N134").Value = -18192.75
$ws.Range("H136").Value = 2502702.5
$ws.Range("I136").Value = 2502702.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7508107.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7505557.5
$ws.Range("N136").Value = ""

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39111676
$ws.Range("I4").Value = 1182122.9
$ws.Range("J4").Value = 178186700
$ws.Range("K4").Value = 3546368.7
$ws.Range("L4").Value = 534560100
$ws.Range("M4").Value = -3546256.7
$ws.Range("N4").Value = -534560324
$ws.Range("H42").Value = 11873.5
$ws.Range("J42").Value = 11873.5
$ws.Range("L42").Value = 35620.5
$ws.Range("N42").Value = -36688.5
$ws.Range("H57").Value = 17153.6
$ws.Range("I57").Value = 5388.5
$ws.Range("K57").Value = 16165.5
$ws.Range("M57").Value = -15606.5
$ws.Range("H96").Value = 12997
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").Value = ""
$ws.Range("H121").Value = 553.3
$ws.Range("J121").Value = 763.25
$ws.Range("L121").Value = 2289.75
$ws.Range("N121").Value = -4909.75
$ws.Range("H131").Value = 2079.2534
$ws.Range("I131").Value = 2066.8462
$ws.Range("J131").Value = 2081.8547
$ws.Range("K131").Value = 6200.5386
$ws.Range("L131").Value = 6245.5641
$ws.Range("M131").Value = -1160.5386
$ws.Range("N131").Value = -16325.5641

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 50495
$ws.Range("J94").Value = 50495
$ws.Range("L94").Value = 50495
$ws.Range("N94").Value = -51847
$ws.Range("H132").Value = 4968.61
$ws.Range("I132").Value = 4523.517
$ws.Range("J132").Value = 6044.25
$ws.Range("K132").Value = 13570.551
$ws.Range("L132").Value = 18132.75
$ws.Range("M132").Value = -11040.551
$ws.Range("N132").Value = -23192.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 514.3333
$ws.Range("I55").Value = 509.33334
$ws.Range("J55").Value = 526.8333
$ws.Range("K55").Value = 509.33334
$ws.Range("L55").Value = 526.8333
$ws.Range("M55").Value = -336.33334
$ws.Range("N55").Value = -872.8333
$ws.Range("H61").Value = 3666.5
$ws.Range("I61").Value = 4274.75
$ws.Range("K61").Value = 4274.75
$ws.Range("M61").Value = -4072.75
$ws.Range("H113").Value = 3666.5
$ws.Range("I113").Value = 4274.75
$ws.Range("K113").Value = 4274.75
$ws.Range("M113").Value = -2104.75
$ws.Range("H122").Value = 4773.2
$ws.Range("I122").Value = 4216.5
$ws.Range("K122").Value = 12649.5
$ws.Range("M122").Value = -10199.5
$ws.Range("H132").Value = 4465.6484
$ws.Range("I132").Value = 3751.309
$ws.Range("J132").Value = 6533.4736
$ws.Range("K132").Value = 11253.927
$ws.Range("L132").Value = 19600.4208
$ws.Range("M132").Value = -8723.927
$ws.Range("N132").Value = -24660.4208
$ws.Range("H141").Value = 102497.5
$ws.Range("J141").Value = 102497.5
$ws.Range("L141").Value = 102497.5
$ws.Range("N141").Value = -112857.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1068.3334
$ws.Range("I113").Value = 1001
$ws.Range("K113").Value = 3003
$ws.Range("M113").Value = -833
$ws.Range("H132").Value = 42742628
$ws.Range("I132").Value = 6174906.5
$ws.Range("K132").Value = 18524719.5
$ws.Range("M132").Value = -18522189.5
$ws.Range("H136").Value = 9145.135
$ws.Range("I136").Value = 3874.6
$ws.Range("J136").Value = 9968.656000000001
$ws.Range("K136").Value = 11623.8
$ws.Range("L136").Value = 29905.968
$ws.Range("M136").Value = -9073.799999999999
$ws.Range("N136").Value = -35005.968
